$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '28.454.60'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = "'" + '1.823.87'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'" + '314.40'
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = "'" + '0.5102'
$ws.Range("E7").Value = '  -4.45%  '
$ws.Range("D8").Value = "'" + '0.3929'
$ws.Range("E8").Value = '  -2.86%  '
$ws.Range("D9").Value = "'" + '0.07698'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").Value = "'" + '41.92'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").Value = "'" + '1.109'
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").Value = "'" + '6.269'
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = "'" + '7.524'
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").Value = "'" + '1.822.00'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = "'" + '93.11'
$ws.Range("E17").Value = '  +4.32%  '
$ws.Range("D18").Value = "'" + '0.00001112'
$ws.Range("E18").Value = '  +3.79%  '
$ws.Range("D19").Value = "'" + '0.06644'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("D20").Value = "'" + '17.75'
$ws.Range("E20").Value = '  +1.10%  '
$ws.Range("D21").Value = "'" + '1.001'
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = "'" + '6.107'
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").Value = "'" + '28.489.42'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = "'" + '11.26'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  +4.91%  '
$ws.Range("D26").Value = "'" + '21.42'
$ws.Range("E26").Value = '  +4.25%  '
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = "'" + '2.034.41'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = "'" + '155.67'
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").Value = "'" + '2.411'
$ws.Range("E29").Value = '  -2.64%  '
$ws.Range("D30").Value = "'" + '124.83'
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = "'" + '0.1100'
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'" + '1.110'
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("D33").Value = "'" + '5.683'
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").Value = "'" + '3.655'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").Value = "'" + '0.07066'
$ws.Range("E35").Value = '  -1.24%  '
$ws.Range("D36").Value = "'" + '0.2215'
$ws.Range("E36").Value = '  -2.29%  '
$ws.Range("D37").Value = "'" + '0.02329'
$ws.Range("E37").Value = '  -0.57%  '
$ws.Range("D38").Value = "'" + '5.182'
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").Value = "'" + '8.779'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").Value = "'" + '0.6277'
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("D41").Value = "'" + '11.21'
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").Value = "'" + '13.40'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").Value = "'" + '3.732'
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("D47").Value = "'" + '0.5894'
$ws.Range("E47").Value = '  +0.78%  '
$ws.Range("E48").Value = '  -1.37%  '
$ws.Range("D49").Value = "'" + '1.985'
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").Value = "'" + '1.194'
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("D51").Value = "'" + '0.06905'
$ws.Range("E51").Value = '  +0.12%  '
